#
# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains a new (blank) column so the
# schedule can carry a "Late" / "Outstanding" pair that is now split by
# an extra spacer column. Concretely: insert one blank column before the
# existing column N (14) on the "Repayment Schedule" sheet, which pushes
# the old N/O columns (headed "Late" / "heading") out to O/P and extends
# the used range from A1:P14 to A1:Q14. The sheet also becomes the active
# tab with a fresh selection.

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column immediately before column N (14th column),
# shifting the "Late"/"heading"/"Outstanding" columns one place to the
# right (N->O, O->P, and a brand new blank P->Q is created for the last
# column's data).
$schedule.Columns.Item(14).Insert()

# Make "Repayment Schedule" the active sheet/tab and move the selection
# to J19, matching the saved view state.
$schedule.Activate()
$schedule.Range("J19").Select()
